$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date serial for every
# data row (rows 2..535). The whole column is bumped from 45202 to 45203.
for ($r = 2; $r -le 535; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

# Rows 533-535 also had their "Beteckning" (A) and "Area (ha)" (G) values
# rotate: the old row-533 record moves down to row 535, and rows 534/535
# shift up into 533/534.
$ws.Cells.Item(533, 1).Value = "A 46922-2023"
$ws.Cells.Item(533, 7).Value = 1.2

$ws.Cells.Item(534, 1).Value = "A 47019-2023"
$ws.Cells.Item(534, 7).Value = 2.2

$ws.Cells.Item(535, 1).Value = "A 47027-2023"
$ws.Cells.Item(535, 7).Value = 7.7
